$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "Rol se encarga de los requisitos de la aplicación" -> split the run into
#    "Rol " + "se encarga de los requisitos de la aplicación"
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Rol se encarga de los requisitos de la aplicación", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($rng.Find.Found) {
    $splitPos = $rng.Start + 4
    $sub = $d.Range($rng.Start, $splitPos)
    # Toggling Bold on then off forces Word to materialize a distinct run
    # boundary at this point without altering the visible formatting.
    $sub.Font.Bold = $true
    $sub.Font.Bold = $false
}

# ---------------------------------------------------------------------------
# 2) Collapse "/Documentación/" + "Gestión" + bookmark(_GoBack) +
#    " de la configuración de software" + ".doc" into a single run:
#    "/Documentación/Gestion_Configuracion.doc"
#    (this also removes the old _GoBack bookmark that lived in that text)
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("/Documentación/Gestión de la configuración de software.doc", $true, $false, $false, $false, $false, $true, 1, $false, "/Documentación/Gestion_Configuracion.doc", 2)

# ---------------------------------------------------------------------------
# 3) "Desarrollo de interfaces web" -> split into
#    "Desarrollo de interfaces " + "web"
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Desarrollo de interfaces web", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($rng.Find.Found) {
    $part1 = "Desarrollo de interfaces "
    $splitPos = $rng.Start + $part1.Length
    $sub = $d.Range($rng.Start, $splitPos)
    $sub.Font.Bold = $true
    $sub.Font.Bold = $false
}

# ---------------------------------------------------------------------------
# 4) "/Diseño de Interfaces/Diseño web/Interfaz.doc" -> split into
#    "/Diseño de Interfaces/Diseño " + "web" + "/Interfaz.doc"
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("/Diseño de Interfaces/Diseño web/Interfaz.doc", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($rng.Find.Found) {
    $part1 = "/Diseño de Interfaces/Diseño "
    $part2 = "web"
    $p1End = $rng.Start + $part1.Length
    $p2End = $p1End + $part2.Length

    $sub1 = $d.Range($rng.Start, $p1End)
    $sub1.Font.Bold = $true
    $sub1.Font.Bold = $false

    $sub2 = $d.Range($p1End, $p2End)
    $sub2.Font.Bold = $true
    $sub2.Font.Bold = $false
}

# ---------------------------------------------------------------------------
# 5) Move the _GoBack bookmark to the final (last) paragraph of the document.
# ---------------------------------------------------------------------------
$paras = $d.Paragraphs
$lastPara = $paras.Last
$lastRng = $lastPara.Range
$lastRng.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:bookmarkStart w:id="1" w:name="_GoBack"/><w:bookmarkEnd w:id="1"/></w:p>')

Write-Output "Edits applied"
